$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Production (MW)" values for rows 2-41 (Timestamp shifts +3 days for all data rows)
$newB = @{
    2 = 135; 3 = 145; 4 = 164; 5 = 169; 6 = 164; 7 = 153; 8 = 144; 9 = 150; 10 = 166;
    11 = 166; 12 = 165; 13 = 161; 14 = 160; 15 = 153; 16 = 153; 17 = 144; 18 = 134; 19 = 126;
    20 = 117; 21 = 120; 22 = 122; 23 = 125; 24 = 132; 25 = 136; 26 = 118; 27 = 97; 28 = 89;
    29 = 81; 30 = 58; 31 = 44; 32 = 30; 33 = 31; 34 = 10; 35 = 9; 36 = 8; 37 = 10; 38 = 11;
    39 = 13; 40 = 10; 41 = 5
}

for ($r = 2; $r -le 97; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $oldVal = [double]$aCell.Value2
    $aCell.Value2 = $oldVal + 3

    if ($newB.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $newB[$r]
    }
}
